$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.759.58"
$ws.Range("E2").Value = "  -4.70%  "

$ws.Range("D3").Value = "2.452.73"
$ws.Range("E3").Value = "  -6.15%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").Value = "'543.53"
$ws.Range("E5").Value = "  -5.53%  "

$ws.Range("D6").Value = "'145.57"
$ws.Range("E6").Value = "  -7.07%  "

$ws.Range("E7").Value = "  -0.06%  "

$ws.Range("D8").Value = "'0.609"
$ws.Range("E8").Value = "  -2.49%  "

$ws.Range("D9").Value = "2.444.55"
$ws.Range("E9").Value = "  -6.37%  "

$ws.Range("D10").Value = "'0.106"
$ws.Range("E10").Value = "  -10.40%  "

$ws.Range("D11").Value = "'0.154"
$ws.Range("E11").Value = "  -1.91%  "

$ws.Range("D12").Value = "'5.33"
$ws.Range("E12").Value = "  -8.77%  "

$ws.Range("D13").Value = "'0.354"
$ws.Range("E13").Value = "  -7.05%  "

$ws.Range("D14").Value = "'25.89"
$ws.Range("E14").Value = "  -8.36%  "

$ws.Range("D15").Value = "2.887.17"
$ws.Range("E15").Value = "  -6.40%  "

$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").Value = "'0.0000162"
$ws.Range("E16").Value = "  -9.82%  "

$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "60.645.89"
$ws.Range("E17").Value = "  -4.71%  "

$ws.Range("D18").Value = "2.437.54"
$ws.Range("E18").Value = "  -6.63%  "

$ws.Range("D19").Value = "'11.07"
$ws.Range("E19").Value = "  -7.92%  "

$ws.Range("D20").Value = "'6.95"
$ws.Range("E20").Value = "  -8.62%  "

$ws.Range("D21").Value = "'4.18"
$ws.Range("E21").Value = "  -7.88%  "

$ws.Range("D22").Value = "'318.20"
$ws.Range("E22").Value = "  -7.45%  "

$ws.Range("E23").Value = "  -0.03%  "

$ws.Range("D24").Value = "'62.97"
$ws.Range("E24").Value = "  -6.74%  "

$ws.Range("D25").Value = "'1.74"
$ws.Range("E25").Value = "  -4.73%  "

$ws.Range("D26").Value = "2.589.47"
$ws.Range("E26").Value = "  -5.21%  "

$ws.Range("B27").Value = "PEPE"
$ws.Range("C27").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D27").Value = "0.0₃0969"
$ws.Range("E27").Value = "  -11.45%  "

$ws.Range("B28").Value = "Binance-PegBSC-USD"
$ws.Range("C28").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D28").Value = "'0.998"
$ws.Range("E28").Value = "  -0.24%  "

$ws.Range("D29").Value = "'1.49"
$ws.Range("E29").Value = "  -5.70%  "

$ws.Range("B30").Value = "Bittensor"
$ws.Range("C30").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D30").Value = "'531.98"
$ws.Range("E30").Value = "  -10.48%  "

$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").Value = "'8.31"
$ws.Range("E31").Value = "  -9.40%  "

$ws.Range("D32").Value = "'7.63"
$ws.Range("E32").Value = "  -3.57%  "

$ws.Range("D33").Value = "'0.147"
$ws.Range("E33").Value = "  -9.17%  "

$ws.Range("D34").Value = "'1.88"
$ws.Range("E34").Value = "  -8.79%  "

$ws.Range("D35").Value = "'1.56"
$ws.Range("E35").Value = "  -10.41%  "

$ws.Range("D36").Value = "'5.81"
$ws.Range("E36").Value = "  -11.80%  "

$ws.Range("B37").Value = "NEARProtocol"
$ws.Range("C37").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D37").Value = "'4.86"
$ws.Range("E37").Value = "  -9.85%  "

$ws.Range("B38").Value = "FirstDigitalUSD"
$ws.Range("C38").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D38").Value = "'0.998"
$ws.Range("E38").Value = "  -0.01%  "

$ws.Range("D39").Value = "'0.374"
$ws.Range("E39").Value = "  -7.07%  "

$ws.Range("D40").Value = "'18.33"
$ws.Range("E40").Value = "  -7.11%  "

$ws.Range("D41").Value = "'143.86"
$ws.Range("E41").Value = "  -6.65%  "

$ws.Range("E42").Value = "  -0.09%  "

$ws.Range("D43").Value = "'1.69"
$ws.Range("E43").Value = "  -9.52%  "

$ws.Range("D44").Value = "'39.92"
$ws.Range("E44").Value = "  -3.50%  "

$ws.Range("D45").Value = "'2.32"
$ws.Range("E45").Value = "  -9.23%  "

$ws.Range("D46").Value = "'146.47"
$ws.Range("E46").Value = "  -7.04%  "

$ws.Range("D47").Value = "'3.57"
$ws.Range("E47").Value = "  -8.59%  "

$ws.Range("D48").Value = "'20.88"
$ws.Range("E48").Value = "  -12.42%  "

$ws.Range("D49").Value = "'0.0530"
$ws.Range("E49").Value = "  -10.25%  "

$ws.Range("D50").Value = "'0.0940"
$ws.Range("E50").Value = "  -6.25%  "

$ws.Range("D51").Value = "'0.579"
$ws.Range("E51").Value = "  -8.35%  "

Write-Host "Applied cryptos update."
